# Timesheet for Week 16 - roll the timesheet forward to the week of 20/04
# (Sun 20/04 - Sat 26/04), updating the "Week of:" date and the seven
# day-of-week row labels to match. This mirrors the author's commit:
# "Made some dates corrections on some weekly Timesheets".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# "Week of:" date field (was 28/04/2014, now 20/04/2014 - the Sunday that
# starts the new reporting week).
$ws.Range("G8").Value = 41749

# Day-of-week column labels for the seven timesheet rows (11-17).
$ws.Range("A11").Value = "Sun 20/04"
$ws.Range("A12").Value = "Mon 21/04"
$ws.Range("A13").Value = "Tue 22/04"
$ws.Range("A14").Value = "Wed 23/04"
$ws.Range("A15").Value = "Thur 24/04"
$ws.Range("A16").Value = "Fri   25/04"
$ws.Range("A17").Value = "Sat  26/04"

# Leave the active selection where the user's edits ended up.
$ws.Range("A17").Select()
